# All db values now parsed from excel, except for recipe_steps
#
# 1. ingredients!D1 header renamed "weight" -> "unit_weight"
# 2. ingredients!A2 typo fix "all purpose flower" -> "all purpose flour"
# 3. ingredients: new row 10 "derp ingredient" / "derp" / "million"
# 4. View/selection state: "steps" becomes the active sheet/tab, each sheet
#    gets its own final selection, matching the author's final click-through.

$wb = $excel.ActiveWorkbook

# --- ingredients sheet -------------------------------------------------
# (written in the order the new shared strings first appear in the target
# workbook: "all purpose flour", "derp ingredient", "derp", "million",
# "unit_weight")
$ingredients = $wb.Worksheets.Item("ingredients")
$ingredients.Range("A2").Value = "all purpose flour"
$ingredients.Range("A10").Value = "derp ingredient"
$ingredients.Range("B10").Value = "derp"
$ingredients.Range("C10").Value = "million"
$ingredients.Range("D1").Value = "unit_weight"

# --- view/selection state -----------------------------------------------
$units = $wb.Worksheets.Item("units")
$units.Activate()
$units.Range("B17").Select()

$ingredients.Activate()
$ingredients.Range("A25").Select()

$step_types = $wb.Worksheets.Item("step_types")
$step_types.Activate()

$recipes = $wb.Worksheets.Item("recipes")
$recipes.Activate()
$recipes.Range("D8").Select()

$recipe_steps = $wb.Worksheets.Item("recipe_steps")
$recipe_steps.Activate()
$recipe_steps.Range("C24").Select()

# "steps" is the sheet left active/selected in the final saved state.
$steps = $wb.Worksheets.Item("steps")
$steps.Activate()
$steps.Range("B14").Select()
